$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.683.91'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.575.95'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.00%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.30'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.46%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.92%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.577.40'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.82%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.56%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.355'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.16'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.79'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.043.80'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.550.62'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.576.20'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.43'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -6.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.74'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '350.83'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.23'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.60'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.31%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.35'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.96%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.86'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -9.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.709.75'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0990'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '531.34'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.17'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.45%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.70%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.132'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.52%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.46'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.20%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.79'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.82%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.77'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.359'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.34'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.13'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.77'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.86%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.43'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₆0286'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.28%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '149.17'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.566'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.72'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.72'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.06%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.88%  '
